$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.55503691331002
$ws.Range("C2").Value = 6.964138495425454
$ws.Range("D2").Value = 7.408760566910734
$ws.Range("F2").Value = 42.76063026919137
$ws.Range("G2").Value = 51.02241540517826
$ws.Range("H2").Value = 19.89556599191068
$ws.Range("J2").Value = 10.79502197558032
$ws.Range("K2").Value = 13.40991361142136
$ws.Range("L2").Value = 11.56925190974325
$ws.Range("N2").Value = 22.44695025191375
$ws.Range("B3").Value = 17.40014836480524
$ws.Range("C3").Value = 6.910647425013416
$ws.Range("D3").Value = 7.400194578092284
$ws.Range("F3").Value = 42.8083268271109
$ws.Range("G3").Value = 51.05561208351045
$ws.Range("H3").Value = 19.93910083686117
$ws.Range("J3").Value = 10.81593845101952
$ws.Range("K3").Value = 13.30158030741861
$ws.Range("L3").Value = 11.5722954281427
$ws.Range("N3").Value = 22.50982070730788
$ws.Range("B4").Value = 17.30828053950183
$ws.Range("C4").Value = 6.876989603667965
$ws.Range("D4").Value = 7.395897216053584
$ws.Range("F4").Value = 42.84589880027156
$ws.Range("G4").Value = 51.08757839454919
$ws.Range("H4").Value = 19.96876376930925
$ws.Range("J4").Value = 10.82968070217524
$ws.Range("K4").Value = 13.2374713277506
$ws.Range("L4").Value = 11.57576052964543
$ws.Range("N4").Value = 22.55030940353797
$ws.Range("B5").Value = 17.27169511313688
$ws.Range("C5").Value = 6.863072993149544
$ws.Range("D5").Value = 7.394389597609566
$ws.Range("F5").Value = 42.8632914745974
$ws.Range("G5").Value = 51.10351411892683
$ws.Range("H5").Value = 19.98158870000274
$ws.Range("J5").Value = 10.83550741092233
$ws.Range("K5").Value = 13.21197726846942
$ws.Range("L5").Value = 11.57757493157949
$ws.Range("N5").Value = 22.56728429321468
$ws.Range("B6").Value = 17.2656726402192
$ws.Range("C6").Value = 6.860750139912335
$ws.Range("D6").Value = 7.394154016701249
$ws.Range("F6").Value = 42.86630518500321
$ws.Range("G6").Value = 51.10633581266634
$ws.Range("H6").Value = 19.98376277549151
$ws.Range("J6").Value = 10.83648863277669
$ws.Range("K6").Value = 13.20778282035649
$ws.Range("L6").Value = 11.5779005371783
$ws.Range("N6").Value = 22.57013170928801
$ws.Range("B7").Value = 17.30778363892502
$ws.Range("C7").Value = 6.876802727092243
$ws.Range("D7").Value = 7.39587589540159
$ws.Range("F7").Value = 42.84612493776351
$ws.Range("G7").Value = 51.08778153581296
$ws.Range("H7").Value = 19.96893374701999
$ws.Range("J7").Value = 10.82975836498509
$ws.Range("K7").Value = 13.23712491896482
$ws.Range("L7").Value = 11.57578336908144
$ws.Range("N7").Value = 22.55053640635026
$ws.Range("B8").Value = 17.50098341040245
$ws.Range("C8").Value = 6.945866080560662
$ws.Range("D8").Value = 7.405608328847517
$ws.Range("F8").Value = 42.77535534194074
$ws.Range("G8").Value = 51.03145572064444
$ws.Range("H8").Value = 19.90996801299477
$ws.Range("J8").Value = 10.80204756924167
$ws.Range("K8").Value = 13.3720760191097
$ws.Range("L8").Value = 11.56997057323996
$ws.Range("N8").Value = 22.46823725726175
$ws.Range("B9").Value = 17.90376370382466
$ws.Range("C9").Value = 7.074671446445501
$ws.Range("D9").Value = 7.432251958898312
$ws.Range("F9").Value = 42.70238319873661
$ws.Range("G9").Value = 51.01301820802362
$ws.Range("H9").Value = 19.81761706755328
$ws.Range("J9").Value = 10.75482313076483
$ws.Range("K9").Value = 13.65465870727267
$ws.Range("L9").Value = 11.571194291128
$ws.Range("N9").Value = 22.32175891346313
$ws.Range("B10").Value = 18.2118433459051
$ws.Range("C10").Value = 7.165045905643391
$ws.Range("D10").Value = 7.456326710581764
$ws.Range("F10").Value = 42.68893973621605
$ws.Range("G10").Value = 51.05564104675904
$ws.Range("H10").Value = 19.76397703691375
$ws.Range("J10").Value = 10.7244376860978
$ws.Range("K10").Value = 13.87159264337249
$ws.Range("L10").Value = 11.57972470116754
$ws.Range("N10").Value = 22.22315340611465
$ws.Range("B11").Value = 18.35410407365573
$ws.Range("C11").Value = 7.205189695992138
$ws.Range("D11").Value = 7.468230610659608
$ws.Range("F11").Value = 42.69154395635798
$ws.Range("G11").Value = 51.08721437471446
$ws.Range("H11").Value = 19.74266337986253
$ws.Range("J11").Value = 10.71154461898718
$ws.Range("K11").Value = 13.97194748337582
$ws.Range("L11").Value = 11.58524659856522
$ws.Range("N11").Value = 22.18023567694113
$ws.Range("B12").Value = 18.4082331819557
$ws.Range("C12").Value = 7.220248105785132
$ws.Range("D12").Value = 7.472872872604346
$ws.Range("F12").Value = 42.69378245367218
$ws.Range("G12").Value = 51.10091859717075
$ws.Range("H12").Value = 19.73503658610456
$ws.Range("J12").Value = 10.70679554854645
$ws.Range("K12").Value = 14.01015883651547
$ws.Range("L12").Value = 11.58757205096248
$ws.Range("N12").Value = 22.16426140021638
$ws.Range("B13").Value = 18.39656476385455
$ws.Range("C13").Value = 7.217011437886542
$ws.Range("D13").Value = 7.471867135461334
$ws.Range("F13").Value = 42.69324468771763
$ws.Range("G13").Value = 51.09788948006398
$ws.Range("H13").Value = 19.73665939213213
$ws.Range("J13").Value = 10.7078124255835
$ws.Range("K13").Value = 14.00192053189059
$ws.Range("L13").Value = 11.58706082360156
$ws.Range("N13").Value = 22.16768940998656
$ws.Range("B14").Value = 18.35855239911935
$ws.Range("C14").Value = 7.20643145230295
$ws.Range("D14").Value = 7.468609850689806
$ws.Range("F14").Value = 42.69170303168752
$ws.Range("G14").Value = 51.08830682309531
$ws.Range("H14").Value = 19.74202701476536
$ws.Range("J14").Value = 10.71115124174423
$ws.Range("K14").Value = 13.97508714110365
$ws.Range("L14").Value = 11.58543322619935
$ws.Range("N14").Value = 22.17891590304468
$ws.Range("B15").Value = 18.3353009719776
$ws.Range("C15").Value = 7.199932143749112
$ws.Range("D15").Value = 7.466632113926215
$ws.Range("F15").Value = 42.69092175544068
$ws.Range("G15").Value = 51.08266467845528
$ws.Range("H15").Value = 19.74537269641668
$ws.Range("J15").Value = 10.71321370500205
$ws.Range("K15").Value = 13.95867722781778
$ws.Range("L15").Value = 11.58446675595464
$ws.Range("N15").Value = 22.18582859403623
$ws.Range("B16").Value = 18.20258488539504
$ws.Range("C16").Value = 7.162402594296259
$ws.Range("D16").Value = 7.455567719632628
$ws.Range("F16").Value = 42.68894489160044
$ws.Range("G16").Value = 51.05382252512134
$ws.Range("H16").Value = 19.76543208501504
$ws.Range("J16").Value = 10.72529894682794
$ws.Range("K16").Value = 13.86506518209087
$ws.Range("L16").Value = 11.57939672702105
$ws.Range("N16").Value = 22.22599710014711
$ws.Range("B17").Value = 18.12167662828056
$ws.Range("C17").Value = 7.139128650024223
$ws.Range("D17").Value = 7.449022269572954
$ws.Range("F17").Value = 42.68996460844599
$ws.Range("G17").Value = 51.03924707764406
$ws.Range("H17").Value = 19.77852884072141
$ws.Range("J17").Value = 10.73295061628984
$ws.Range("K17").Value = 13.8080430419785
$ws.Range("L17").Value = 11.57670572545305
$ws.Range("N17").Value = 22.25113493100191
$ws.Range("B18").Value = 18.07534205577191
$ws.Range("C18").Value = 7.125651212771963
$ws.Range("D18").Value = 7.445347274435401
$ws.Range("F18").Value = 42.69137200249132
$ws.Range("G18").Value = 51.03201081261047
$ws.Range("H18").Value = 19.78635234308762
$ws.Range("J18").Value = 10.73743915940892
$ws.Range("K18").Value = 13.77540452685002
$ws.Range("L18").Value = 11.57531252952542
$ws.Range("N18").Value = 22.26577605788373
$ws.Range("B19").Value = 18.05969000760951
$ws.Range("C19").Value = 7.121072497700258
$ws.Range("D19").Value = 7.444118474613197
$ws.Range("F19").Value = 42.6919895441352
$ws.Range("G19").Value = 51.02975785687679
$ws.Range("H19").Value = 19.78905114792967
$ws.Range("J19").Value = 10.73897394421256
$ws.Range("K19").Value = 13.76438195741659
$ws.Range("L19").Value = 11.57486741616525
$ws.Range("N19").Value = 22.27076466929884
$ws.Range("B20").Value = 18.13026890564116
$ws.Range("C20").Value = 7.141615635914802
$ws.Range("D20").Value = 7.449709769663445
$ws.Range("F20").Value = 42.68977110871837
$ws.Range("G20").Value = 51.04067995727951
$ws.Range("H20").Value = 19.77710459120869
$ws.Range("J20").Value = 10.73212702942884
$ws.Range("K20").Value = 13.81409690313554
$ws.Range("L20").Value = 11.5769761981387
$ws.Range("N20").Value = 22.24844008405307
$ws.Range("B21").Value = 18.36971091439888
$ws.Range("C21").Value = 7.209542967474526
$ws.Range("D21").Value = 7.469562963282364
$ws.Range("F21").Value = 42.69212188138161
$ws.Range("G21").Value = 51.09107407815296
$ws.Range("H21").Value = 19.74043835600033
$ws.Range("J21").Value = 10.71016693715805
$ws.Range("K21").Value = 13.98296333035811
$ws.Range("L21").Value = 11.5859049419955
$ws.Range("N21").Value = 22.17561088016156
$ws.Range("B22").Value = 18.5276830085142
$ws.Range("C22").Value = 7.253101100189155
$ws.Range("D22").Value = 7.483321023927765
$ws.Range("F22").Value = 42.70095658035189
$ws.Range("G22").Value = 51.13419583945335
$ws.Range("H22").Value = 19.71906421489232
$ws.Range("J22").Value = 10.69659127406122
$ws.Range("K22").Value = 14.09453192973126
$ws.Range("L22").Value = 11.59310600030418
$ws.Range("N22").Value = 22.12963118700005
$ws.Range("B23").Value = 18.44324980333062
$ws.Range("C23").Value = 7.229931119447663
$ws.Range("D23").Value = 7.475907281588233
$ws.Range("F23").Value = 42.69557423883795
$ws.Range("G23").Value = 51.11025059648095
$ws.Range("H23").Value = 19.73023500083623
$ws.Range("J23").Value = 10.70376593854752
$ws.Range("K23").Value = 14.03488573747843
$ws.Range("L23").Value = 11.58913826618174
$ws.Range("N23").Value = 22.15402365888813
$ws.Range("B24").Value = 18.1263837707915
$ws.Range("C24").Value = 7.140491571255275
$ws.Range("D24").Value = 7.449398676373341
$ws.Range("F24").Value = 42.68985603196737
$ws.Range("G24").Value = 51.04002859056615
$ws.Range("H24").Value = 19.77774757859589
$ws.Range("J24").Value = 10.73249909424381
$ws.Range("K24").Value = 13.8113595007341
$ws.Range("L24").Value = 11.5768534379493
$ws.Range("N24").Value = 22.24965783535664
$ws.Range("B25").Value = 17.79248700049386
$ws.Range("C25").Value = 7.040558472740519
$ws.Range("D25").Value = 7.424245660522115
$ws.Range("F25").Value = 42.71507102914425
$ws.Range("G25").Value = 51.00814180708233
$ws.Range("H25").Value = 19.84010607257466
$ws.Range("J25").Value = 10.76683964330297
$ws.Range("K25").Value = 13.57645581693575
$ws.Range("L25").Value = 11.56951801484705
$ws.Range("N25").Value = 22.3597970231451
